$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the sheet; this shifts all existing rows
# (and the merged cells, styles, etc.) down by one.
$ws.Rows.Item(1).Insert()

# Put the new run label in the new row 1 (column C, matching the value columns below)
$ws.Range("C1").Value = "Run(e) 1"

# Update the active selection to match the post-edit workbook state
$ws.Range("E2").Select()
